$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add %DiffH, %DiffD, %DiffA in columns AF, AG, AH,
# reusing the same header style as the existing Diff columns (AC:AE).
$ws.Range("AC1").Copy()
$ws.Range("AF1:AH1").PasteSpecial(-4122)
$ws.Range("AF1").Value = "%DiffH"
$ws.Range("AG1").Value = "%DiffD"
$ws.Range("AH1").Value = "%DiffA"

# Data rows 2..17: %DiffX = DiffX / YtrueX * 100
for ($r = 2; $r -le 17; $r++) {
    $diffH = $ws.Range("AC$r").Value()
    $diffD = $ws.Range("AD$r").Value()
    $diffA = $ws.Range("AE$r").Value()
    $trueH = $ws.Range("Z$r").Value()
    $trueD = $ws.Range("AA$r").Value()
    $trueA = $ws.Range("AB$r").Value()

    $ws.Range("AF$r").Value = ($diffH / $trueH) * 100
    $ws.Range("AG$r").Value = ($diffD / $trueD) * 100
    $ws.Range("AH$r").Value = ($diffA / $trueA) * 100
}
